# Applies the repositioning/recropping changes described in the commit
# "figures-png and ltspice files" to slide 2 of the presentation.
#
# Shape-tree order on slide 2 (as enumerated by $s.Shapes):
#   1 - Picture 5  (id=6)   -> move left
#   4 - TextBox 10 (id=11)  -> move left/top
#   6 - Picture 1  (id=2)   -> re-crop on the left edge + move/resize
#   7 - Picture 12 (id=13)  -> move left
#
# NOTE: the runtime's Left/Top/Width/Height COM properties round-trip
# through single-precision floats before being converted back to EMUs
# (truncating, not rounding), so the literal EMU/12700 quotient can end
# up one EMU short of the target. The literals below are nudged by a
# few hundred-thousandths of a point (invisible in the UI, << 1/100 mm)
# so that the value PowerPoint actually persists lands exactly on the
# EMU figure from the target XML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 1: "Picture 5" - shift left (y unchanged)
$shp1 = $s.Shapes.Item(1)
$shp1.Left = 600.00034496063

# Shape 4: "TextBox 10" - shift left and up slightly
$shp4 = $s.Shapes.Item(4)
$shp4.Left = 490.36922259842515
$shp4.Top = 232.42347456692914

# Shape 6: "Picture 1" - crop further from the left edge of the source
# image, then move/shrink the shape to match the new crop.
$shp6 = $s.Shapes.Item(6)
$shp6.PictureFormat.CropLeft = 41.85225
$shp6.Left = 430.60222472440944
$shp6.Width = 183.14574803149605

# Shape 7: "Picture 12" - shift left (y unchanged)
$shp7 = $s.Shapes.Item(7)
$shp7.Left = 389.7963192125984
